$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.138156
$ws.Range("H2").Value = 0.414468
$ws.Range("I2").Value = 0.0003010053794496939
$ws.Range("J2").Value = 0.0003010053794496939
$ws.Range("M2").Value = 0.8596446666666667
$ws.Range("N2").Value = 2.578934
$ws.Range("O2").Value = 0.05286426382906832
$ws.Range("P2").Value = 0.05286426382906832
$ws.Range("Q2").Value = 0.118765068568
$ws.Range("R2").Value = 1.068885617112
$ws.Range("S2").Value = 0.00001591242779319743
$ws.Range("T2").Value = 0.00001591242779319744
$ws.Range("G3").Value = 0.138156
$ws.Range("H3").Value = 0.414468
$ws.Range("I3").Value = 0.0003010053794496939
$ws.Range("J3").Value = 0.0003010053794496939
$ws.Range("O3").Value = 0.6417658132713033
$ws.Range("P3").Value = 0.6417658132713032
$ws.Range("Q3").Value = 1.441793667348
$ws.Range("R3").Value = 12.976143006132
$ws.Range("S3").Value = 0.00019317496214157
$ws.Range("T3").Value = 0.00019317496214157
$ws.Range("G4").Value = 0.138156
$ws.Range("H4").Value = 0.414468
$ws.Range("I4").Value = 0.0003010053794496939
$ws.Range("J4").Value = 0.0003010053794496939
$ws.Range("O4").Value = 0.3053699228996285
$ws.Range("P4").Value = 0.3053699228996284
$ws.Range("Q4").Value = 0.686045301776
$ws.Range("R4").Value = 6.174407715984
$ws.Range("S4").Value = 0.00009191798951492644
$ws.Range("T4").Value = 0.00009191798951492644
$ws.Range("I5").Value = 0.9878623917146768
$ws.Range("J5").Value = 0.9878623917146769
$ws.Range("M5").Value = 0.8596446666666667
$ws.Range("N5").Value = 2.578934
$ws.Range("O5").Value = 0.05286426382906832
$ws.Range("P5").Value = 0.05286426382906832
$ws.Range("Q5").Value = 389.7722522508938
$ws.Range("R5").Value = 3507.950270258044
$ws.Range("S5").Value = 0.05222261810241911
$ws.Range("T5").Value = 0.05222261810241911
$ws.Range("I6").Value = 0.9878623917146768
$ws.Range("J6").Value = 0.9878623917146769
$ws.Range("O6").Value = 0.6417658132713033
$ws.Range("P6").Value = 0.6417658132713032
$ws.Range("Q6").Value = 4731.788326140225
$ws.Range("R6").Value = 42586.09493526202
$ws.Range("S6").Value = 0.6339763112189044
$ws.Range("T6").Value = 0.6339763112189044
$ws.Range("I7").Value = 0.9878623917146768
$ws.Range("J7").Value = 0.9878623917146769
$ws.Range("O7").Value = 0.3053699228996285
$ws.Range("P7").Value = 0.3053699228996284
$ws.Range("S7").Value = 0.3016634623933535
$ws.Range("T7").Value = 0.3016634623933535
$ws.Range("G8").Value = 5.432785666666668
$ws.Range("I8").Value = 0.01183660290587349
$ws.Range("J8").Value = 0.01183660290587349
$ws.Range("M8").Value = 0.8596446666666667
$ws.Range("N8").Value = 2.578934
$ws.Range("O8").Value = 0.05286426382906832
$ws.Range("P8").Value = 0.05286426382906832
$ws.Range("Q8").Value = 4.670265223493113
$ws.Range("R8").Value = 42.03238701143801
$ws.Range("S8").Value = 0.000625733298856013
$ws.Range("T8").Value = 0.000625733298856013
$ws.Range("G9").Value = 5.432785666666668
$ws.Range("I9").Value = 0.01183660290587349
$ws.Range("J9").Value = 0.01183660290587349
$ws.Range("O9").Value = 0.6417658132713033
$ws.Range("P9").Value = 0.6417658132713032
$ws.Range("S9").Value = 0.007596327090257375
$ws.Range("T9").Value = 0.007596327090257374
$ws.Range("G10").Value = 5.432785666666668
$ws.Range("I10").Value = 0.01183660290587349
$ws.Range("J10").Value = 0.01183660290587349
$ws.Range("O10").Value = 0.3053699228996285
$ws.Range("P10").Value = 0.3053699228996284
$ws.Range("Q10").Value = 26.97774314667956
$ws.Range("S10").Value = 0.003614542516760107
$ws.Range("T10").Value = 0.003614542516760107
